# Adds two new weekly records (Coliflor, Vega Monumental Concepción) for
# 2022-01-06 (serial 44567) just above the existing 2021-12-09 record,
# pushing every subsequent row down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 113-114; everything currently at/after row 113
# (old rows 113..194) shifts down to 115..196 automatically.
$ws.Range("A113:A114").EntireRow.Insert()

# ---- New row 113: Coliflor, Primera, 2022-01-06, Región Metropolitana ----
$ws.Cells.Item(113, 1).Value = 11
$ws.Cells.Item(113, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(113, 3).Value = "Bíobío"
$ws.Cells.Item(113, 4).Value = 44567
$ws.Cells.Item(113, 5).Value = 8
$ws.Cells.Item(113, 6).Value = 100112008
$ws.Cells.Item(113, 7).Value = "Coliflor"
$ws.Cells.Item(113, 8).Value = "Sin especificar"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 2000
$ws.Cells.Item(113, 11).Value = 700
$ws.Cells.Item(113, 12).Value = 800
$ws.Cells.Item(113, 13).Value = 750
$ws.Cells.Item(113, 14).Value = "`$/unidad"
$ws.Cells.Item(113, 15).Value = "Región Metropolitana"
$ws.Cells.Item(113, 16).Value = 750
$ws.Cells.Item(113, 17).Value = 1
$ws.Cells.Item(113, 18).Value = "Hortaliza"

# ---- New row 114: Coliflor, Segunda, 2022-01-06, Región Metropolitana ----
$ws.Cells.Item(114, 1).Value = 11
$ws.Cells.Item(114, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(114, 3).Value = "Bíobío"
$ws.Cells.Item(114, 4).Value = 44567
$ws.Cells.Item(114, 5).Value = 8
$ws.Cells.Item(114, 6).Value = 100112008
$ws.Cells.Item(114, 7).Value = "Coliflor"
$ws.Cells.Item(114, 8).Value = "Sin especificar"
$ws.Cells.Item(114, 9).Value = "Segunda"
$ws.Cells.Item(114, 10).Value = 1000
$ws.Cells.Item(114, 11).Value = 600
$ws.Cells.Item(114, 12).Value = 600
$ws.Cells.Item(114, 13).Value = 600
$ws.Cells.Item(114, 14).Value = "`$/unidad"
$ws.Cells.Item(114, 15).Value = "Región Metropolitana"
$ws.Cells.Item(114, 16).Value = 600
$ws.Cells.Item(114, 17).Value = 1
$ws.Cells.Item(114, 18).Value = "Hortaliza"
